$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Results data block (rows 14-20) - entered first so later AVERAGE formulas
# have real data to compute against, and so the "0.00" / "0" number formats
# get minted in the same order the final styles.xml expects (2-decimal
# format first, then integer format, then percent format).
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "1."
$ws.Range("B14").Value = 8.83
$ws.Range("C14").Value = 978

$ws.Range("A15").Value = "2"
$ws.Range("B15").Value = 6.98
$ws.Range("C15").Value = 773

$ws.Range("A16").Value = "3"
$ws.Range("B16").Value = 11.02
$ws.Range("C16").Value = 1162

$ws.Range("A17").Value = "4"
$ws.Range("B17").Value = 11.49
$ws.Range("C17").Value = 1221

$ws.Range("A18").Value = "5"
$ws.Range("B18").Value = 65.34
$ws.Range("C18").Value = 7589

$ws.Range("A19").Value = "6"
$ws.Range("B19").Value = 9.94
$ws.Range("C19").Value = 1091

$ws.Range("A20").Value = "7"
$ws.Range("B20").Value = 49.17
$ws.Range("C20").Value = 5543

# 2-decimal number format for the "time" column -> mints numFmtId 2 first
$ws.Range("B14:B20").NumberFormat = "0.00"
# plain integer number format for the "generation" column -> mints numFmtId 1
$ws.Range("C14:C20").NumberFormat = "0"

# ---------------------------------------------------------------------------
# CONFIG block (rows 4-8): re-apply values + integer number format
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 100
$ws.Range("B5").Value = 10000
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 20
$ws.Range("B8").Value = 60

$ws.Range("B4").NumberFormat = "0"
$ws.Range("B5").NumberFormat = "0"
$ws.Range("B6").NumberFormat = "0"
$ws.Range("B7").NumberFormat = "0"
$ws.Range("B8").NumberFormat = "0"

$ws.Range("C4").NumberFormat = "0"
$ws.Range("C5").NumberFormat = "0"
$ws.Range("C6").NumberFormat = "0"
$ws.Range("C7").NumberFormat = "0"
$ws.Range("C8").NumberFormat = "0"
$ws.Range("C9").NumberFormat = "0"
$ws.Range("C10").NumberFormat = "0"
$ws.Range("C11").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Row 9 used to hold "Vyriešených" in A9 - that label moved down to A10, so
# clear A9 entirely.
# ---------------------------------------------------------------------------
$ws.Range("A9").ClearContents()

# ---------------------------------------------------------------------------
# Row 10: "Vyriešených" + solved ratio, percent-formatted
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Vyriešených"
$ws.Range("B10").Value = 0.7
$ws.Range("B10").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# Row 11: average solving time (renamed label + formula, 2-decimal format)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Priemerný čas vyriešených (s)"
$ws.Range("B11").Formula = "=AVERAGE(B14:B20)"
$ws.Range("B11").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Row 12: average final generation (new row, integer format)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Priemerná konečná generácia"
$ws.Range("B12").Formula = "=AVERAGE(C14:C20)"
$ws.Range("B12").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Row 13: table headers, shifted down from row 13/14 before -> still row 13
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Dáta"
$ws.Range("B13").Value = "ČAS"
$ws.Range("C13").Value = "GENERÁCIA RIEŠENIA"

# ---------------------------------------------------------------------------
# Rows 21-23: unfinished runs
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "8"
$ws.Range("B21").Value = "NEDOKONČIL"

$ws.Range("A22").Value = "9"
$ws.Range("B22").Value = "NEDOKONČIL"

$ws.Range("A23").Value = "10"
$ws.Range("B23").Value = "NEDOKONČIL"

# ---------------------------------------------------------------------------
# Selection, matching the saved cursor position in the authored file
# ---------------------------------------------------------------------------
$null = $ws.Range("I18").Select()
